# Generate Report for Handback
# Updates the handoff/handback timestamp cells for the
# "ca6ecd48-c673-4fea-9792-4020c15d8bc0" row across the Overview,
# zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 3
$wsOverview.Range("G3").Value = "2016-08-26 20:47:45"

# zh-cn sheet: row 3 is the ca6ecd48-... entry
#   H = Correspond Handoff Datetime
#   K = Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-26 20:47:41"
$wsZhCn.Range("K3").Value = "2016-08-26 20:47:57"

# de-de sheet: row 3 is the ca6ecd48-... entry
#   H = Correspond Handoff Datetime
#   K = Correspond Handback DateTime
$wsDeDe.Range("H3").Value = "2016-08-26 20:47:45"
$wsDeDe.Range("K3").Value = "2016-08-26 20:48:10"
